$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "Tue Sep 26 21:18:20 EDT 2023",
    "Tue Sep 26 21:18:30 EDT 2023",
    "Tue Sep 26 21:18:40 EDT 2023",
    "Tue Sep 26 21:18:50 EDT 2023",
    "Tue Sep 26 21:18:59 EDT 2023",
    "Tue Sep 26 21:19:09 EDT 2023",
    "Tue Sep 26 21:19:19 EDT 2023",
    "Tue Sep 26 21:19:29 EDT 2023",
    "Tue Sep 26 21:19:39 EDT 2023",
    "Tue Sep 26 21:19:49 EDT 2023",
    "Tue Sep 26 21:19:59 EDT 2023",
    "Tue Sep 26 21:20:08 EDT 2023",
    "Tue Sep 26 21:20:18 EDT 2023",
    "Tue Sep 26 21:20:28 EDT 2023",
    "Tue Sep 26 21:20:38 EDT 2023",
    "Tue Sep 26 21:20:48 EDT 2023"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value = $dates[$i]
}
